$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet,
#    using the existing "2021-Q4" sheet as a formatting template so that the
#    header row / index column keep the same bold+border style.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet    = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other data sheets (a brand-new sheet
# otherwise defaults to Excel's standard 0.7"/0.75" margins).
$newSheet.PageSetup.LeftMargin   = 54
$newSheet.PageSetup.RightMargin  = 54
$newSheet.PageSetup.TopMargin    = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy the template layout (formats + values) then overwrite the values.
$templateSheet.Range("A1:H3").Copy($newSheet.Range("A1"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'501307"
$newSheet.Cells.Item(2,3).Value = "银河中证沪港深高股息指数（LOF）A"
$newSheet.Cells.Item(2,4).Value = "'0.19"
$newSheet.Cells.Item(2,5).Value = "'91.35"
$newSheet.Cells.Item(2,6).Value = "'1.46"
$newSheet.Cells.Item(2,7).Value = "'0.0028"
$newSheet.Cells.Item(2,8).Value = 9

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'501308"
$newSheet.Cells.Item(3,3).Value = "银河中证沪港深高股息指数（LOF）C"
$newSheet.Cells.Item(3,4).Value = "'0.01"
$newSheet.Cells.Item(3,5).Value = "'91.35"
$newSheet.Cells.Item(3,6).Value = "'1.46"
$newSheet.Cells.Item(3,7).Value = "'0.0001"
$newSheet.Cells.Item(3,8).Value = 9

# ---------------------------------------------------------------------------
# 2. Insert a new row 2 at the top of the "总计" sheet for the 2022-Q1 totals,
#    pushing the previous rows down, then renumber the index column (A).
#    NOTE: re-fetch the sheet by name now that the sheet collection has
#    changed (the previously captured $totalSheet reference no longer
#    points at the "总计" tab after the insertion above).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2").EntireRow.Insert()

# re-apply the row formatting (bold/bordered index cell) to the new row
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2"))

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 2
$totalSheet.Cells.Item(2,4).Value = 0

for ($i = 1; $i -le 5; $i++) {
    $totalSheet.Cells.Item($i + 2, 1).Value = $i
}
